$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 1-9: update column A text (existing rows) ---
$ws.Range("A1").Value = "Patient education, counseling"
$ws.Range("A2").Value = "Lab tests are needed to evaluate patient status today."
$ws.Range("A3").Value = "Tests, procedures, other laboratory studies : "
$ws.Range("A4").Value = "Specialist referral(s) or consults"
$ws.Range("A5").Value = " “The client shows greater success with activates involving physical cueing.”"
$ws.Range("A6").Value = "“Continue with POC.” (plan of care)"
$ws.Range("A7").Value = "."
$ws.Range("A8").Value = "."
$ws.Range("A9").Value = "."

# --- Row 10: Periarthritis bullet list (A) + food/sodium advice (B) ---
$a10 = @"
• Periarthritis shoulder – Pain and limitation of movements 
• Cheiroarthropathy – Stiffness and limited mobility of fingers 
• Neuroarthropathy • Painless unilateral swelling of foot and ankle 
• Spinal hyperostosis – Mild back pain with preservation of back movements
"@
$ws.Range("A10").Value = $a10
$b10 = @"
Buy fresh food often.
Avoid  to “fast” foods, frozen dinners, and canned foods.
Use spices, herbs, and sodium-free seasonings in place of salt. 
Check for sodium on the Nutrition Facts label of food packages. 
Rinse canned vegetables, beans, meats etc  before eating.
Look for food labels with words like sodium free or salt free
"@
$ws.Range("B10").Value = $b10
$ws.Rows.Item(10).RowHeight = 64.15
$ws.Range("A10:B10").WrapText = $true

# --- Row 11: weight loss / fat reduction (A) + sodium diet target (B) ---
$a11 = @"
Weight loss of 5%-10%.  
Reduction in fat intake < 30% of calories.  
Reduction in saturated fat intake < 10% of calories.
"@
$ws.Range("A11").Value = $a11
$ws.Range("B11").Value = "Your diet should contain less than 2,300 milligrams of sodium each day."
$ws.Rows.Item(11).RowHeight = 44
$ws.Range("A11:B11").WrapText = $true

# --- Row 12: physical activity (A) + citation link (B) ---
$a12 = @"
Increase in physical activity levels. 
The brisk walking should last for at least 30 minutes 
and should be undertaken  at least three times a week.
"@
$ws.Range("A12").Value = $a12
$ws.Range("B12").Value = "https://owl.purdue.edu/owl/subject_specific_writing/healthcare_writing/soap_notes/soap_note_tips.html"
$ws.Rows.Item(12).RowHeight = 44
$ws.Range("A12:B12").WrapText = $true

# --- Row 13: foot inspection reminder ---
$a13 = @"
Regular inspection and examination of the foot at risk.
Appropriate footwear.
"@
$ws.Range("A13").Value = $a13
$ws.Rows.Item(13).RowHeight = 29.85
$ws.Range("A13").WrapText = $true

# --- Rows 14-19: blank padding rows (extend used range to B19) ---
for ($r = 14; $r -le 18; $r++) {
  $ws.Rows.Item($r).RowHeight = 12.8
}
$ws.Rows.Item(19).RowHeight = 12.8
$ws.Range("B19").NumberFormat = "General"

# --- Selection matches the commit focus cell ---
$ws.Range("A6").Select()
